$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.298.22'
$ws.Range("E2").Value = '  +1.16%  '
$ws.Range("D3").Value = '2.482.75'
$ws.Range("E3").Value = '  +3.10%  '
$ws.Range("E4").Value = '  -0.28%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '577.53'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.71%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.70'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.71%  '
$ws.Range("E7").Value = '  +0.25%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("D9").Value = '2.484.38'
$ws.Range("E9").Value = '  +1.94%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.111'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.28%  '
$ws.Range("E11").Value = '  +1.79%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.25'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.47%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.353'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.48%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.58'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.11%  '
$ws.Range("E15").Value = '  +1.12%  '
$ws.Range("D16").Value = '2.936.26'
$ws.Range("E16").Value = '  +1.82%  '
$ws.Range("D17").Value = '63.186.37'
$ws.Range("E17").Value = '  +0.91%  '
$ws.Range("D18").Value = '2.484.21'
$ws.Range("E18").Value = '  +2.41%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '8.13'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.83%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.04'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.93%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '330.14'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.75%  '
$ws.Range("E22").Value = '  +9.13%  '
$ws.Range("E23").Value = '  +0.03%  '
$ws.Range("E24").Value = '  +0.31%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '66.29'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.07%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.86'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +15.76%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '663.06'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +6.28%  '
$ws.Range("D28").Value = '0.0₃0996'
$ws.Range("E28").Value = '  +1.53%  '
$ws.Range("D29").Value = '2.608.27'
$ws.Range("E29").Value = '  +1.84%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -9.50%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.48'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +5.23%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.08'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.24%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.86'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.97%  '
$ws.Range("E34").Value = '  -3.52%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.55'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.23%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.999'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.32%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.79'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.85%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.48'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.83%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.372'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.47%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.78'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.93%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '150.37'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.02%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.69'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.70%  '
$ws.Range("E43").Value = '  +0.52%  '
$ws.Range("E44").Value = '  +0.03%  '
$ws.Range("E45").Value = '  -43.50%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '155.62'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +7.84%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '15.24'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.33%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.62'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.69%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '20.47'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.21%  '
$ws.Range("E50").Value = '  +1.92%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0515'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.16%  '
